$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H9:H35").NumberFormat = '_-* #,##0.00_-;\-* #,##0.00_-;_-* "-"??_-;_-@_-'
